# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 6386dc03-....ee27304... file row in both the zh-cn and
# de-de language sheets to reflect the new report generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 02:49:26"
$wsZhCn.Range("H3").Value = "2016-03-24 02:49:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 02:49:30"
$wsDeDe.Range("H3").Value = "2016-03-24 02:49:57"
